$d = $word.ActiveDocument

$d.Content.Find.Execute("819÷7=117, 0", $true, $true, $false, $false, $false, $true, 1, $false, "661÷2=330, 1", 2) | Out-Null
$d.Content.Find.Execute("445÷9=49, 4", $true, $true, $false, $false, $false, $true, 1, $false, "250÷6=41, 4", 2) | Out-Null
$d.Content.Find.Execute("576÷4=144, 0", $true, $true, $false, $false, $false, $true, 1, $false, "408÷6=68, 0", 2) | Out-Null
$d.Content.Find.Execute("785÷2=392, 1", $true, $true, $false, $false, $false, $true, 1, $false, "476÷2=238, 0", 2) | Out-Null
$d.Content.Find.Execute("975÷5=195, 0", $true, $true, $false, $false, $false, $true, 1, $false, "155÷7=22, 1", 2) | Out-Null
$d.Content.Find.Execute("434÷3=144, 2", $true, $true, $false, $false, $false, $true, 1, $false, "786÷9=87, 3", 2) | Out-Null
$d.Content.Find.Execute("855÷7=122, 1", $true, $true, $false, $false, $false, $true, 1, $false, "300÷4=75, 0", 2) | Out-Null
$d.Content.Find.Execute("475÷9=52, 7", $true, $true, $false, $false, $false, $true, 1, $false, "139÷4=34, 3", 2) | Out-Null
$d.Content.Find.Execute("872÷7=124, 4", $true, $true, $false, $false, $false, $true, 1, $false, "467÷5=93, 2", 2) | Out-Null
$d.Content.Find.Execute("946÷9=105, 1", $true, $true, $false, $false, $false, $true, 1, $false, "285÷4=71, 1", 2) | Out-Null
$d.Content.Find.Execute("498÷6=83, 0", $true, $true, $false, $false, $false, $true, 1, $false, "306÷9=34, 0", 2) | Out-Null
$d.Content.Find.Execute("500÷9=55, 5", $true, $true, $false, $false, $false, $true, 1, $false, "188÷6=31, 2", 2) | Out-Null
$d.Content.Find.Execute("185÷3=61, 2", $true, $true, $false, $false, $false, $true, 1, $false, "468÷9=52, 0", 2) | Out-Null
$d.Content.Find.Execute("186÷2=93, 0", $true, $true, $false, $false, $false, $true, 1, $false, "172÷7=24, 4", 2) | Out-Null
$d.Content.Find.Execute("672÷9=74, 6", $true, $true, $false, $false, $false, $true, 1, $false, "912÷7=130, 2", 2) | Out-Null
$d.Content.Find.Execute("437÷7=62, 3", $true, $true, $false, $false, $false, $true, 1, $false, "374÷5=74, 4", 2) | Out-Null
$d.Content.Find.Execute("290÷5=58, 0", $true, $true, $false, $false, $false, $true, 1, $false, "821÷6=136, 5", 2) | Out-Null
$d.Content.Find.Execute("471÷7=67, 2", $true, $true, $false, $false, $false, $true, 1, $false, "663÷5=132, 3", 2) | Out-Null
$d.Content.Find.Execute("633÷6=105, 3", $true, $true, $false, $false, $false, $true, 1, $false, "238÷9=26, 4", 2) | Out-Null
$d.Content.Find.Execute("154÷6=25, 4", $true, $true, $false, $false, $false, $true, 1, $false, "691÷8=86, 3", 2) | Out-Null
$d.Content.Find.Execute("750÷2=375, 0", $true, $true, $false, $false, $false, $true, 1, $false, "279÷7=39, 6", 2) | Out-Null
$d.Content.Find.Execute("824÷9=91, 5", $true, $true, $false, $false, $false, $true, 1, $false, "478÷4=119, 2", 2) | Out-Null
$d.Content.Find.Execute("844÷2=422, 0", $true, $true, $false, $false, $false, $true, 1, $false, "373÷7=53, 2", 2) | Out-Null
$d.Content.Find.Execute("360÷7=51, 3", $true, $true, $false, $false, $false, $true, 1, $false, "434÷7=62, 0", 2) | Out-Null
$d.Content.Find.Execute("321÷3=107, 0", $true, $true, $false, $false, $false, $true, 1, $false, "583÷4=145, 3", 2) | Out-Null

$d.Save()
